$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Patient identification block ---
# Patient name (merged A6:F7)
$ws.Range("A6").Value = "AJÙ  LINARES  ROSA  AMANDA"
# Expediente No. (merged G6:H7) - numeric-looking, force text with quote prefix
$ws.Range("G6").Value = "'5348"

# Fecha de Nacimiento - date-looking text, force text with quote prefix
$ws.Range("A9").Value = "'1999-05-20"
# Edad - numeric-looking, force text with quote prefix
$ws.Range("D9").Value = "'18"
# Lugar de Nacimiento
$ws.Range("E9").Value = "GUATEMALTECO"

# Ocupacion
$ws.Range("C11").Value = "AMA DE CASA"
# Documento de Identificacion - numeric-looking, force text with quote prefix
$ws.Range("G11").Value = "'2997248400101"

# --- Emergency contact block ---
# Nombre
$ws.Range("A13").Value = "AMANDA LINARES"
# Direccion
$ws.Range("E13").Value = "MANZANA D LOTE 23 EL MESQUITAL"
# Telefono - numeric-looking, force text with quote prefix
$ws.Range("G13").Value = "'32081430"

# --- Attention info block ---
# Hora
$ws.Range("D14").Value = "Hora: 11:40"
# Area de urgencia
$ws.Range("E14").Value = "Area de urgencia: GINECOLOGIA"
# Fecha de la asistencia medica - date-looking text, force text with quote prefix
$ws.Range("A15").Value = "'2017-10-16"

# Tipo de Consulta value cell - cleared out
$ws.Range("D16").Value = ""
